$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "locacalizacion" -> "localizacion"
$ws.Range("B1").Value = "localizacion"

# Move selection to the edited cell
$ws.Range("B1").Select()
